# Adding the changes we made on may 9th
#
# The underlying sensor data window was re-sampled/shifted: for every row,
# the elapsed-time column (A, 0,100,200,... ms) stays fixed to the row's
# position, but the ax/ay/az/gx/gy/gz readings (columns C-H) are replaced.
# Net effect: the data grows from 20 rows (A1:H21) to 30 rows (A1:H31) -
# 7 new samples are effectively inserted in front (using data that used to
# be further down the series) and 3 new samples are appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, "walkingToRunning", 11.77978897094727, -1.277233123779297, 9.730603218078612, -0.4386182389586079, -0.4680532807871848, 2.060500017882589),
    @(100, "walkingToRunning", 1.413437724113464, -7.171976566314697, 6.755977630615234, 0.6352223770569996, -0.8322304098058368, -1.049409866333014),
    @(200, "walkingToRunning", -4.193170547485352, -4.375148296356201, 1.568653106689453, 0.1475016089714127, -0.1921065187431543, -0.1445163721927848),
    @(300, "walkingToRunning", -1.307081580162048, -2.807691097259521, -1.606552600860596, -1.213809009395561, -0.2346711329983395, 1.28885372863315),
    @(400, "walkingToRunning", -1.827142477035522, 0.1487736701965332, 4.701539993286133, -0.7734762763008961, -0.09618946024911929, 2.230298755737722),
    @(500, "walkingToRunning", -1.866428852081299, 4.58729076385498, 1.570873260498047, -0.814207781389884, 0.3517244944944581, 1.427822031950592),
    @(600, "walkingToRunning", 1.589986324310303, -0.6990594863891602, 7.147370338439941, -0.4451152733739872, 0.1002355693144553, 1.227771341498127),
    @(700, "walkingToRunning", 11.77978897094727, -1.277233123779297, 9.730603218078612, -0.4386182389586079, -0.4680532807871848, 2.060500017882589),
    @(800, "walkingToRunning", 4.897948265075684, 1.876962661743164, 4.303222179412842, -1.567447010333154, -0.1997661534507753, 1.985428574121546),
    @(900, "walkingToRunning", 17.21874618530273, -8.779304504394531, -3.521855354309082, -0.6623021185133424, 1.233293045912556, 0.3937855432481396),
    @(1000, "walkingToRunning", 4.876998901367188, 3.794233560562134, -0.4165830612182617, 1.065984005429924, -0.5608204790782889, -2.052491584707599),
    @(1100, "walkingToRunning", 0.1654682159423828, -41.74227905273438, 9.125425338745115, 1.985182776668951, -1.370061922194394, -2.959458767460082),
    @(1200, "walkingToRunning", 16.16831207275391, -3.528035402297974, 29.5200309753418, 3.586229367909687, 3.220478899587801, -3.701172142464491),
    @(1300, "walkingToRunning", -0.2412894368171692, 4.501626968383789, 4.947979927062988, -2.049730680920693, -0.3309398520415856, 0.1050501378054163),
    @(1400, "walkingToRunning", -14.82570934295654, -18.40331840515137, -2.339614391326904, -6.480390057345936, -0.290497637037169, 3.6122629448847),
    @(1500, "walkingToRunning", -2.732851982116699, 11.58825302124023, -3.568616628646851, -3.746594315252966, 5.263269504314723, 2.846692088896868),
    @(1600, "walkingToRunning", -8.033671379089355, -1.333880424499511, -8.920849800109863, 0.02857317537218518, 1.524775265437166, 0.3577531473285853),
    @(1700, "walkingToRunning", 13.30204772949219, -27.35689926147461, 14.50175476074219, 4.081077046200642, -3.786933429350092, -1.819117558486549),
    @(1800, "walkingToRunning", -3.511280059814453, 11.19524669647217, 6.339620590209961, 5.00775977923788, -1.871817700754, -3.02260407941595),
    @(1900, "walkingToRunning", 2.829581499099731, -56.68034744262695, 7.749513149261475, -0.1797645285649305, 0.1183369633207527, -0.06624912005396677),
    @(2000, "walkingToRunning", -12.45886611938477, 26.94960594177246, -9.934419631958008, -7.62000698970656, -6.875881371537381, 5.026607214496812),
    @(2100, "walkingToRunning", -15.69838333129883, 4.704256057739258, -6.9078369140625, -3.768456836642528, -15.13744506739117, 5.25606083264814),
    @(2200, "walkingToRunning", -22.81047821044922, -8.102127075195312, 3.835212230682373, 4.952609313926157, -11.08443266002066, 3.078852102841224),
    @(2300, "walkingToRunning", 26.74444198608398, 23.38513946533203, 12.87228775024414, 4.517102752240144, 0.1121426887317769, 1.888300155020016),
    @(2400, "walkingToRunning", 15.44845008850098, 10.10699462890625, 11.63222694396973, 3.696813958550474, 4.523056784559633, 1.050867166615984),
    @(2500, "walkingToRunning", -1.591778755187989, -18.45157623291016, -9.295619010925291, -3.479645563261123, 1.03966360149645, 1.598440124298719),
    @(2600, "walkingToRunning", 14.80349922180176, 22.9969539642334, -1.046570301055908, -4.971113100874858, 5.150270287760582, 0.4502635256288277),
    @(2700, "walkingToRunning", -12.98141479492188, -5.830618858337402, -7.131386756896973, -2.691276774793723, 7.824305781253101, -2.586581079487865),
    @(2800, "walkingToRunning", 15.02213287353516, -47.13114547729492, 8.131996154785156, 1.383527442585964, 2.926498572838676, -5.056850963437613),
    @(2900, "walkingToRunning", -2.31139087677002, 14.0579719543457, 5.640069961547852, 3.296718087898288, -3.505119464119212, -4.962326313638459)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
